$wb = $excel.ActiveWorkbook

# --- Sheet "Bico": clear the "Validado com sucesso!..." messages in column H
#     for the rows where the reconciliation message was removed (bug fix),
#     leaving rows 12 and 13 untouched.
$bico = $wb.Worksheets.Item("Bico")
$bicoRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 14, 15)
foreach ($r in $bicoRows) {
    $bico.Cells.Item($r, 8).Value = ""
}

# --- Sheet "Tanque": replace the "Validado com sucesso!..." messages in
#     column F with the "Divergência..." messages produced by the buggy
#     comparison (SPED value stuck at row 4's closing value, report value
#     built by concatenating the closing value's digits instead of being
#     multiplied).
$tanque = $wb.Worksheets.Item("Tanque")
$tanque.Cells.Item(2, 6).Value = "Divergência entre o SPED(5389,63) e o relatório(1411587,00)!"
$tanque.Cells.Item(3, 6).Value = "Divergência entre o SPED(5389,63) e o relatório(5025131,00)!"
$tanque.Cells.Item(4, 6).Value = "Divergência entre o SPED(5389,63) e o relatório(5389627,00)!"
$tanque.Cells.Item(5, 6).Value = "Divergência entre o SPED(5389,63) e o relatório(194023,00)!"
